# Weekly update: a new daily price observation for
# "Feria Lagunitas de Puerto Montt - Pepino ensalada" is inserted as the
# new row 241, pushing the existing rows 241-360 down to 242-361.
#
# Excel's native row-insert does exactly this: it shifts all the rows
# below down by one and leaves the freshly inserted row blank (but keeps
# column formatting, e.g. the date style on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 241; rows 241:360 shift down to 242:361,
# and the sheet's dimension grows from R360 to R361 automatically.
$ws.Rows.Item(241).Insert()

# The row that lands on 242 is a duplicate of what used to be row 241 -
# reuse it as a template for the new row's unchanged metadata columns
# (Mercado, Region, Codreg, Categoria, Variedad, Calidad, Precio maximo,
# Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion).
$ws.Range("A242:R242").Copy()
$ws.Range("A241:R241").PasteSpecial()

# Now overwrite the cells that actually hold the new observation's data:
# Fecha, Volumen, Precio minimo, Precio promedio ponderado, Precio $/Kg.
$ws.Cells.Item(241, 4).Value  = 44917   # D241 Fecha
$ws.Cells.Item(241, 10).Value = 120     # J241 Volumen
$ws.Cells.Item(241, 11).Value = 22000   # K241 Precio minimo
$ws.Cells.Item(241, 13).Value = 22500   # M241 Precio promedio ponderado
$ws.Cells.Item(241, 16).Value = 375     # P241 Precio $/Kg
